$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M106").Value = -1967.1292
$ws.Range("K106").Value = 2598.1292
$ws.Range("H106").Value = 2727.4211
$ws.Range("I106").Value = 2598.1292
$ws.Range("N111").ClearContents()
$ws.Range("H111").Value = 126775.375
$ws.Range("M111").Value = -377259.125
$ws.Range("L111").Value = 0
$ws.Range("K111").Value = 380326.125
$ws.Range("J111").Value = 0
$ws.Range("I111").Value = 126775.375
$ws.Range("L112").Value = 5505.1764
$ws.Range("H112").Value = 1835.0588
$ws.Range("J112").Value = 1835.0588
$ws.Range("N112").Value = -7721.1764
$ws.Range("K138").Value = 8991.875100000001
$ws.Range("L138").Value = 22307.946
$ws.Range("J138").Value = 7435.982
$ws.Range("M138").Value = -3851.875100000001
$ws.Range("N138").Value = -32587.946
$ws.Range("H138").Value = 6104.375
$ws.Range("I138").Value = 2997.2917

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K4").Value = 159
$ws.Range("I4").Value = 159
$ws.Range("M4").Value = -43
$ws.Range("J4").Value = 547
$ws.Range("N4").Value = -779
$ws.Range("L4").Value = 547
$ws.Range("H4").Value = 288.33334
$ws.Range("K45").Value = 1764.4783
$ws.Range("H45").Value = 2095.6667
$ws.Range("I45").Value = 1764.4783
$ws.Range("M45").Value = -1387.4783
$ws.Range("M61").Value = -2870.9412
$ws.Range("K61").Value = 3082.9412
$ws.Range("I61").Value = 3082.9412
$ws.Range("H61").Value = 3032.3215
$ws.Range("H74").Value = 2511.5652
$ws.Range("K74").Value = 2491.625
$ws.Range("I74").Value = 2491.625
$ws.Range("M74").Value = -1617.625
$ws.Range("K77").Value = 12458.125
$ws.Range("H77").Value = 2511.5652
$ws.Range("I77").Value = 2491.625
$ws.Range("M77").Value = -8090.125
$ws.Range("L102").Value = 3911
$ws.Range("M102").Value = -889.3076000000001
$ws.Range("I102").Value = 2511.3076
$ws.Range("K102").Value = 2511.3076
$ws.Range("H102").Value = 2611.2856
$ws.Range("N102").Value = -7155
$ws.Range("J102").Value = 3911
$ws.Range("M136").Value = -6698.8236
$ws.Range("K136").Value = 9248.8236
$ws.Range("I136").Value = 3082.9412
$ws.Range("H136").Value = 3032.3215

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 411
$ws.Range("I25").Value = 411
$ws.Range("K25").Value = 411
$ws.Range("M25").Value = -176
$ws.Range("L57").Value = 79492
$ws.Range("N57").Value = -80932
$ws.Range("H57").Value = 79492
$ws.Range("J57").Value = 79492
$ws.Range("H107").Value = 296070.66
$ws.Range("K107").Value = 1812.9667
$ws.Range("M107").Value = 107.0333000000001
$ws.Range("I107").Value = 1812.9667
$ws.Range("L136").Value = 79492
$ws.Range("J136").Value = 79492
$ws.Range("H136").Value = 79492
$ws.Range("N136").Value = -89692

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2875
$ws.Range("N62").Value = -4201.5715
$ws.Range("L62").Value = 2953.5715
$ws.Range("J62").Value = 2953.5715
$ws.Range("N65").Value = -21007.8575
$ws.Range("L65").Value = 14767.8575
$ws.Range("H65").Value = 2875
$ws.Range("J65").Value = 2953.5715
$ws.Range("N107").Value = -4837.5
$ws.Range("H107").Value = 919.38464
$ws.Range("K107").Value = 905.1818
$ws.Range("M107").Value = 1014.8182
$ws.Range("L107").Value = 997.5
$ws.Range("J107").Value = 997.5
$ws.Range("I107").Value = 905.1818
$ws.Range("K132").Value = 3960.15
$ws.Range("H132").Value = 1693.0698
$ws.Range("M132").Value = -1430.15
$ws.Range("I132").Value = 1320.05

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K139").Value = 8180.4288
$ws.Range("I139").Value = 2726.8096
$ws.Range("M139").Value = -3040.4288
$ws.Range("J139").Value = 9068.210999999999
$ws.Range("H139").Value = 5738.975
$ws.Range("L139").Value = 27204.633
$ws.Range("N139").Value = -37484.633

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N63").Value = -46372
$ws.Range("H63").Value = 45000
$ws.Range("J63").Value = 45000
$ws.Range("L63").Value = 45000
$ws.Range("L66").Value = 135000
$ws.Range("N66").Value = -141864
$ws.Range("J66").Value = 45000
$ws.Range("H66").Value = 45000
$ws.Range("I80").Value = 1256361.6
$ws.Range("M80").Value = -1255363.6
$ws.Range("K80").Value = 1256361.6
$ws.Range("H80").Value = 1115787.1
$ws.Range("H83").Value = 1115787.1
$ws.Range("K83").Value = 6281808
$ws.Range("I83").Value = 1256361.6
$ws.Range("M83").Value = -6276816
$ws.Range("L113").Value = 16312.125
$ws.Range("J113").Value = 16312.125
$ws.Range("H113").Value = 634804.25
$ws.Range("N113").Value = -20652.125
$ws.Range("K122").Value = 7883.0868
$ws.Range("I122").Value = 2627.6956
$ws.Range("H122").Value = 2903.0967
$ws.Range("M122").Value = -5433.0868
$ws.Range("J123").Value = 69994.25
$ws.Range("L123").Value = 69994.25
$ws.Range("H123").Value = 69994.25
$ws.Range("N123").Value = -74894.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 54063
$ws.Range("J6").Value = 54063
$ws.Range("N6").Value = -54287
$ws.Range("L6").Value = 54063
$ws.Range("H40").Value = 3404606.2
$ws.Range("M40").Value = -7287164.5
$ws.Range("I40").Value = 7287300.5
$ws.Range("K40").Value = 7287300.5
$ws.Range("L82").Value = 2333
$ws.Range("M82").Value = -1369.8
$ws.Range("N82").Value = -3055
$ws.Range("K82").Value = 1730.8
$ws.Range("H82").Value = 1956.625
$ws.Range("I82").Value = 1730.8
$ws.Range("J82").Value = 2333
$ws.Range("J85").Value = 2333
$ws.Range("N85").Value = -4829
$ws.Range("K85").Value = 1730.8
$ws.Range("I85").Value = 1730.8
$ws.Range("M85").Value = -482.8
$ws.Range("L85").Value = 2333
$ws.Range("H85").Value = 1956.625
$ws.Range("L92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("H96").Value = 80852
$ws.Range("N96").Value = -86344
$ws.Range("L96").Value = 80852
$ws.Range("J96").Value = 80852
$ws.Range("N109").Value = -64774
$ws.Range("L109").Value = 62000
$ws.Range("H109").Value = 62000
$ws.Range("J109").Value = 62000

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N107").Value = -5562
$ws.Range("H107").Value = 878.6818
$ws.Range("K107").Value = 2978.8125
$ws.Range("M107").Value = -1058.8125
$ws.Range("L107").Value = 1722
$ws.Range("J107").Value = 574
$ws.Range("I107").Value = 992.9375
$ws.Range("N111").Value = -143242
$ws.Range("H111").Value = 135062
$ws.Range("L111").Value = 135062
$ws.Range("J111").Value = 135062
$ws.Range("K122").Value = 157898028
$ws.Range("I122").Value = 52632676
$ws.Range("H122").Value = 28574134
$ws.Range("M122").Value = -157895578
